$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "User Story #1"
$ws.Range("C3").Value = 13

$ws.Range("B4").Value = "Update GUI to have perimter as option"
$ws.Range("D4").Value = "John"
$ws.Range("E4").Value = "Complete"
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0

$ws.Range("B5").Value = "Multiple Options in GUI"
$ws.Range("D5").Value = "John "
$ws.Range("E5").Value = "Complete"
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1

$ws.Range("B6").Value = "Implement sending emails"
$ws.Range("D6").Value = "John"
$ws.Range("E6").Value = "In Progress"
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

$ws.Range("B7").Value = "Task"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""

$ws.Range("B8").Value = "User Story #2"
$ws.Range("C8").Value = 15

$ws.Range("B9").Value = "Learn Google ApI"
$ws.Range("D9").Value = "Emily"
$ws.Range("E9").Value = "Complete"
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1

$ws.Range("B10").Value = "Attach Google API to code"
$ws.Range("D10").Value = "Emily"
$ws.Range("E10").Value = "In Progress"
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

$ws.Range("B11").Value = "Attach the perimter option to GUI"
$ws.Range("D11").Value = "Emily"
$ws.Range("E11").Value = "Complete"
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0

$ws.Range("B12").Value = "Task"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""

$ws.Range("B13").Value = "User Story #3"
$ws.Range("C13").Value = 12

$ws.Range("B14").Value = "Build UML "
$ws.Range("D14").Value = "Isaac"
$ws.Range("E14").Value = "Complete"
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0

$ws.Range("B15").Value = "Double check the GUI"
$ws.Range("D15").Value = "Isaac"
$ws.Range("E15").Value = "Complete"
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0

$ws.Range("B16").Value = "Edit the UML "
$ws.Range("D16").Value = "Isaac "
$ws.Range("E16").Value = "Complete"
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1

$ws.Range("B17").Value = "Task"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""

$ws.Range("B18").Value = "User Story #4"
$ws.Range("C18").Value = 15

$ws.Range("B19").Value = "Build Gui "
$ws.Range("D19").Value = "Jarod "
$ws.Range("E19").Value = "Complete"
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0

$ws.Range("B20").Value = "Redo the Sprint Backlog"
$ws.Range("D20").Value = "Jarod "
$ws.Range("E20").Value = "In Progress"
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0

$ws.Range("B21").Value = "Merge both GUI from sprints"
$ws.Range("D21").Value = "Jarod"
$ws.Range("E21").Value = "Complete"
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0

$ws.Range("B22").Value = "Task "

$ws.Range("B23").Value = "User Story #5"

$ws.Range("B24").Value = "Task"

$ws.Range("B25").Value = "Task"

$ws.Range("B26").Value = "Task"

$ws.Range("B27").Value = "Task"

$ws.Range("B16").Select()
